# Insert one new data row right before the existing row 837, shifting all
# subsequent rows (837-926) down to (838-927). This matches the diff's
# dimension change from A1:R926 to A1:R927.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(837).Insert()

# Populate the newly inserted row 837 with its data.
$ws.Range("A837").Value = 10
$ws.Range("B837").Value = "Vega Modelo de Temuco"
$ws.Range("C837").Value = "La Araucanía"
$ws.Range("D837").Value = 45194
$ws.Range("E837").Value = 9
$ws.Range("F837").Value = 100112045
$ws.Range("G837").Value = "Zapallo"
$ws.Range("H837").Value = "Camote"
$ws.Range("I837").Value = "1a (guarda)"
$ws.Range("J837").Value = 1000
$ws.Range("K837").Value = 1000
$ws.Range("L837").Value = 1100
$ws.Range("M837").Value = 1050
$ws.Range("N837").Value = "$/kilo (volumen en unidades)"
$ws.Range("O837").Value = "Perú"
$ws.Range("P837").Value = 1050
$ws.Range("Q837").Value = 1
$ws.Range("R837").Value = "Hortaliza"
